# Update the "AMIs - Instance Store Backed" sheet: the AMI ids/names for
# rows 2-8 have been re-pulled from the source and shuffled into a new order.
$wb = $excel.ActiveWorkbook

$amiSheet = $wb.Worksheets.Item("AMIs - Instance Store Backed")

$amiRows = @(
    @("ami-01b2b400e24fcdbe5", "jitsi-2-latest"),
    @("ami-0d1213ebe53bab3a3", "jitsi-2"),
    @("ami-0eb42d6e6d10db5d1", "jitsi-5"),
    @("ami-0ef04d7a197cea5c8", "jitsi-2-official"),
    @("ami-03fe541226ed7a78f", "jitsi-final"),
    @("ami-0c757bca3918f6fa9", "jitsi"),
    @("ami-0f84d671dcbc8a5e9", "jitsi-4")
)

for ($i = 0; $i -lt $amiRows.Length; $i++) {
    $row = $i + 2
    $amiSheet.Range("A$row").Value = $amiRows[$i][0]
    $amiSheet.Range("B$row").Value = $amiRows[$i][1]
}

# Update the "Lambda - Functions" sheet: LastModifiedDaysAgo for the "test"
# function moved from 8 to 9.
$lambdaSheet = $wb.Worksheets.Item("Lambda - Functions")
$lambdaSheet.Range("K3").Value = 9

# Add a new "RDS - Instances" sheet at the end of the workbook containing a
# single informational cell.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$rdsSheet = $wb.Worksheets.Add($null, $lastSheet)
$rdsSheet.Name = "RDS - Instances"
$rdsSheet.Range("A1").Value = "Invalid data format."
